$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("D382")

# Update the confidentiality note date from 2021-04-05 to 2021-04-06
$noteCell = $ws.Range("A58")
$noteText = $noteCell.Value()
$noteCell.Value = $noteText -replace "2021-04-05", "2021-04-06"

# Update Weight (D) and Percent Change (E) values for rows 2-55
$ws.Range("D2").Value = 0.01618326058470936
$ws.Range("E2").Value = 0.006525198938992061
$ws.Range("D3").Value = 0.05130072937076258
$ws.Range("E3").Value = -0.0009018418026918162
$ws.Range("D4").Value = 0.01479816807843185
$ws.Range("E4").Value = 0.007198263821740936
$ws.Range("D5").Value = 0.009574817619152534
$ws.Range("E5").Value = 0.005446333687566529
$ws.Range("D6").Value = 0.01547716866105467
$ws.Range("E6").Value = 0.003562447611064501
$ws.Range("D7").Value = 0.02039211663598169
$ws.Range("E7").Value = 0.0003664345914253797
$ws.Range("D8").Value = 0.004267553956841395
$ws.Range("E8").Value = -0.009075262087310132
$ws.Range("D9").Value = 0.006467580711130532
$ws.Range("E9").Value = -0.002831858407079557
$ws.Range("D10").Value = 0.01386179975600623
$ws.Range("E10").Value = 0.004562533548040815
$ws.Range("D11").Value = 0.008803413941413958
$ws.Range("E11").Value = -0.002745069710323067
$ws.Range("D12").Value = 0.01458735165702184
$ws.Range("E12").Value = 0.002615746795710194
$ws.Range("D13").Value = 0.003128382144860262
$ws.Range("E13").Value = 0.01743151903237283
$ws.Range("D14").Value = 0.006051035444091151
$ws.Range("E14").Value = 0.00144508670520227
$ws.Range("D15").Value = 0.01416565521950457
$ws.Range("E15").Value = -0.007030334591850207
$ws.Range("D16").Value = 0.01040965700740589
$ws.Range("E16").Value = -0.002061855670103196
$ws.Range("D17").Value = 0.02177371143391002
$ws.Range("E17").Value = -0.003942958533219465
$ws.Range("D18").Value = 0.008434946265501564
$ws.Range("E18").Value = -0.0186827105763141
$ws.Range("D19").Value = 0.01662925019657916
$ws.Range("E19").Value = -0.0002447531053051399
$ws.Range("D20").Value = 0.01175802357601695
$ws.Range("E20").Value = -0.01130401860565733
$ws.Range("D21").Value = 0.007226074665325751
$ws.Range("E21").Value = 0.008791924454575106
$ws.Range("D22").Value = 0.01336200903025623
$ws.Range("E22").Value = -0.005597014925373234
$ws.Range("D23").Value = 0.01908508662058828
$ws.Range("E23").Value = 0.003958614484930401
$ws.Range("D24").Value = 0.009609444931809472
$ws.Range("E24").Value = 0.005565684675175264
$ws.Range("D25").Value = 0.021164315247436
$ws.Range("E25").Value = -0.001171875000000044
$ws.Range("D26").Value = 0.01148186360317443
$ws.Range("E26").Value = 0.005931956964233764
$ws.Range("D27").Value = 0.02311285677141868
$ws.Range("E27").Value = -0.01259079903147686
$ws.Range("D28").Value = 0.05604600669340548
$ws.Range("E28").Value = 0.002462271644162017
$ws.Range("D29").Value = 0.02183676557624125
$ws.Range("E29").Value = -0.0175901495162708
$ws.Range("D30").Value = 0.0327067527997486
$ws.Range("E30").Value = -0.01368850865253735
$ws.Range("D31").Value = 0.01650985115247741
$ws.Range("E31").Value = -0.0161761562032352
$ws.Range("D32").Value = 0.01341453825019127
$ws.Range("E32").Value = 0.01918573230048626
$ws.Range("D33").Value = 0.02146321032454823
$ws.Range("E33").Value = -0.002773333333333405
$ws.Range("D34").Value = 0.04233422682823069
$ws.Range("E34").Value = -0.004371417240508935
$ws.Range("D35").Value = 0.01099297936800417
$ws.Range("E35").Value = -0.003471017007983246
$ws.Range("D36").Value = 0.009449917633737979
$ws.Range("E36").Value = -0.007160354249104994
$ws.Range("D37").Value = 0.01198362576453942
$ws.Range("E37").Value = 0.02786069651741285
$ws.Range("D38").Value = 0.007118154464079061
$ws.Range("E38").Value = 0.01701956580005337
$ws.Range("D39").Value = 0.01176947062152338
$ws.Range("E39").Value = -0.002763819095477404
$ws.Range("D40").Value = 0.0176648898414244
$ws.Range("E40").Value = -0.007345926349933252
$ws.Range("D41").Value = 0.01720134809311141
$ws.Range("E41").Value = -0.008710311885361088
$ws.Range("D42").Value = 0.0337220421414713
$ws.Range("E42").Value = 0.005657548584198535
$ws.Range("D43").Value = 0.01122147511525188
$ws.Range("E43").Value = 0.0008646495100677054
$ws.Range("D44").Value = 0.02159116285543116
$ws.Range("E44").Value = -0.001780496712929014
$ws.Range("D45").Value = 0.01397913197244709
$ws.Range("E45").Value = -0.009437438584998348
$ws.Range("D46").Value = 0.008057861508110725
$ws.Range("E46").Value = 0.002130908832617129
$ws.Range("D47").Value = 0.01337644502653377
$ws.Range("E47").Value = -0.01313593770056909
$ws.Range("D48").Value = 0.009543401838707121
$ws.Range("E48").Value = 0.003688377114069885
$ws.Range("D49").Value = 0.01495785436324649
$ws.Range("E49").Value = -0.002946353409351787
$ws.Range("D50").Value = 0.008240759857424503
$ws.Range("E50").Value = -0.01134798294522787
$ws.Range("D51").Value = 0.01119368423255016
$ws.Range("E51").Value = 0.007987910189982461
$ws.Range("D52").Value = 0.008572533393019078
$ws.Range("E52").Value = -0.001917662899343142
$ws.Range("D53").Value = 0.1380796048530767
$ws.Range("E53").Value = -0.0001971220185293943
$ws.Range("D54").Value = 0.04386603230108223
$ws.Range("E54").Value = -0.0003696857670979492
$ws.Range("D55").Value = 0.9999999999999999
$ws.Range("E55").Value = -0.001379618910684499

$ws.Protect("D382")
